$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap A19 and B19 values
$a19 = $ws.Range("A19").Value2
$b19 = $ws.Range("B19").Value2
$ws.Range("A19").Value2 = $b19
$ws.Range("B19").Value2 = $a19

# Set column A width (target OOXML width=27; ColumnWidth setter has a
# constant +5/6 padding baked into the exported <col width> value)
$ws.Columns("A").ColumnWidth = 26.16666666666667

# Scroll/selection state
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("A37").Select()
